$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-22 Thursday" "2026-01-23 Friday"

Replace-Text "781×6=4686" "113×6=678"
Replace-Text "214×8=1712" "742×6=4452"
Replace-Text "842×5=4210" "551×3=1653"
Replace-Text "473×3=1419" "559×4=2236"
Replace-Text "787×8=6296" "373×9=3357"

Replace-Text "916×9=8244" "297×2=594"
Replace-Text "479×8=3832" "837×9=7533"
Replace-Text "179×5=895" "658×8=5264"
Replace-Text "141×7=987" "399×5=1995"
Replace-Text "578×3=1734" "705×8=5640"

Replace-Text "307×8=2456" "555×5=2775"
Replace-Text "512×6=3072" "152×9=1368"
Replace-Text "499×4=1996" "621×6=3726"
Replace-Text "870×4=3480" "249×8=1992"
Replace-Text "903×2=1806" "826×7=5782"

Replace-Text "756×2=1512" "121×8=968"
Replace-Text "749×7=5243" "314×2=628"
Replace-Text "675×4=2700" "784×3=2352"
Replace-Text "610×4=2440" "499×4=1996"
Replace-Text "977×7=6839" "858×5=4290"

Replace-Text "438×2=876" "638×2=1276"
Replace-Text "785×2=1570" "632×4=2528"
Replace-Text "418×5=2090" "285×6=1710"
Replace-Text "714×6=4284" "238×7=1666"
Replace-Text "361×7=2527" "440×4=1760"
